$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.4335615932941437, 0.1406517177820205, -0.8185594081878662),
    @(0.09498954564332961, -0.7519751191139221, -0.1093448773026466),
    @(0.1846340149641037, -1.312596678733826, 0.0687223374843597),
    @(0.6478226184844971, -0.9091202020645142, -0.1838704347610473),
    @(-0.1064432710409164, -0.09178250283002851, 0.0652098655700683),
    @(-0.042302418500185, 0.3572034537792206, 0.1937969923019409),
    @(-0.2768746614456177, 0.2338086664676666, -0.1817324161529541),
    @(0.0734565481543541, 0.1968513280153274, 0.1055269688367843),
    @(-0.1507309973239898, -0.0175623763352632, 0.08170322328805921),
    @(-0.1734857261180877, 0.0183259565383195, 0.0250454749912023)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
